# Append a new match-log row (row 3) to the "Trent Boult " sheet.
# It duplicates the existing row 2 entry (Mumbai Indians vs Chennai Super
# Kings, Abu Dhabi, September 19 2020) — same batsman line, 0 runs off 1
# ball, 0 fours, 0 sixes, strike rate 0.00.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = "3"

$ws.Range("A" + $row).Value = " Abu Dhabi"
$ws.Range("B" + $row).Value = " September 19 2020"
$ws.Range("C" + $row).Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D" + $row).Value = "Mumbai Indians"
$ws.Range("E" + $row).Value = "Chennai Super Kings"
$ws.Range("F" + $row).Value = "Trent Boult "

# The stat columns (totalRuns/totalBalls/total4s/total6s/sr) are stored as
# plain text in this sheet, just like in row 2, even though they look like
# numbers. Writing them as a text formula and then pasting-as-values over
# themselves converts them to literal text cells without leaving behind
# any new/explicit cell style (equivalent to typing the formula, copying,
# and doing Paste Special > Values in the Excel UI).
$statRange = $ws.Range("G" + $row + ":K" + $row)
$ws.Range("G" + $row).Formula = "=""0"""
$ws.Range("H" + $row).Formula = "=""1"""
$ws.Range("I" + $row).Formula = "=""0"""
$ws.Range("J" + $row).Formula = "=""0"""
$ws.Range("K" + $row).Formula = "=""0.00"""

$statRange.Copy()
$statRange.PasteSpecial(-4163)
